# Auto-generated Excel COM-interop script applying the cryptos.xlsx update
# commit: "Updated cryptos list on Sun May 26 17:40:54 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.878.52"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "3.854.08"
$ws.Range("E3").Value = "  +3.01%  "
$ws.Range("E4").Value = "  +0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "601.36"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.09%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "162.63"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.69%  "
$ws.Range("D7").Value = "3.852.33"
$ws.Range("E7").Value = "  +3.05%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("E10").Value = "  -0.53%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "6.32"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("E12").Value = "  -0.30%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "36.78"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -3.06%  "
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").Value = "4.493.72"
$ws.Range("E15").Value = "  +2.90%  "
$ws.Range("D16").Value = "3.845.50"
$ws.Range("E16").Value = "  +2.26%  "
$ws.Range("D17").Value = "69.027.30"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("E18").Value = "  +2.63%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.43"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +4.05%  "
$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.113"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.26%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.14"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.68%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "484.81"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("E24").Value = "  +4.13%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "83.70"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.30%  "
$ws.Range("E26").Value = "  -2.20%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "12.08"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.00%  "
$ws.Range("E28").Value = "  -0.86%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.97"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.73%  "
$ws.Range("E31").Value = "  -1.26%  "
$ws.Range("D32").Value = "4.002.98"
$ws.Range("E33").Value = "  -3.94%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "32.25"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.43%  "
$ws.Range("D35").Value = "3.799.10"
$ws.Range("E36").Value = "  -1.86%  "
$ws.Range("E37").Value = "  +1.52%  "
$ws.Range("E38").Value = "  +3.91%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "5.88"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  -1.72%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "441.97"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.41%  "
$ws.Range("E43").Value = "  -0.22%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "48.47"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("E47").Value = "  -1.04%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "26.61"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +12.91%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "142.77"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.16%  "
$ws.Range("D50").Value = "2.829.63"
$ws.Range("E50").Value = "  +1.66%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0358"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.66%  "
